# Append three blank paragraphs and a final paragraph of text to the end
# of the document, right after the "-> write a program in java." line,
# carrying that paragraph's run/mark formatting (Times New Roman (Body
# CS) complex-script font, sz=16) but with no literal text run on the
# blank ones -- just like Word leaves a paragraph mark with no run when
# you press Enter without typing anything.

$d = $word.ActiveDocument

# Locate the very last paragraph ("-> write a program in java.") and get
# a zero-length range collapsed to the end of the document (right after
# its paragraph mark).
$count = $d.Paragraphs.Count
$lastRange = $d.Paragraphs.Item($count).Range
$insertPoint = $d.Range($lastRange.End, $lastRange.End)

# Build a WordprocessingML package fragment containing the four new
# paragraphs and insert it in one shot via Range.InsertXML so that the
# empty paragraphs come through with only <w:pPr> (no synthesized run).
$rPr = '<w:rPr><w:rFonts w:cs="Times New Roman (Body CS)"/><w:sz w:val="16"/></w:rPr>'
$emptyPara = "<w:p><w:pPr>$rPr</w:pPr></w:p>"
$textPara = "<w:p><w:pPr>$rPr</w:pPr><w:r>$rPr<w:t>This is a change in the file dear</w:t></w:r></w:p>"

$body = $emptyPara + $emptyPara + $emptyPara + $textPara

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       "<w:body>$body</w:body>" +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

[void]$insertPoint.InsertXML($xml)
